$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.017.27"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = "'1.820.05"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = "'310.96"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('D7').Value = "'0.4479"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.02%  '
$ws.Range('D8').Value = "'0.3693"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').Value = "'0.8543"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('D12').Value = "'1.816.35"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').Value = "'6.639"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('D14').Value = "'92.38"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.79%  '
$ws.Range('D15').Value = "'0.07095"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = "'5.322"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = "'0.000008786"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').Value = "'14.95"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('D21').Value = "'26.982.03"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('D22').Value = "'5.158"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').Value = "'10.91"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('D25').Value = "'151.59"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.12%  '
$ws.Range('D26').Value = "'2.228"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.06%  '
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('D28').Value = "'5.239"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').Value = "'116.37"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = "'0.08844"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.58%  '
$ws.Range('D31').Value = "'1.181"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').Value = "'0.7515"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').Value = "'2.966"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.68%  '
$ws.Range('D34').Value = "'4.442"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('D36').Value = "'1.094"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('D37').Value = "'0.01966"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').Value = "'0.05231"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = "'0.5316"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.51%  '
$ws.Range('D40').Value = "'7.159"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +0.55%  '
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('D43').Value = "'0.5257"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.91%  '
$ws.Range('D44').Value = "'8.503"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.50%  '
$ws.Range('D45').Value = "'10.68"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.21%  '
$ws.Range('D46').Value = "'1.968"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.62%  '
$ws.Range('D47').Value = "'105.52"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.66%  '
$ws.Range('D48').Value = "'0.9999"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('D49').Value = "'1.666"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('D50').Value = "'0.06380"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('E51').Value = '  +0.32%  '
